$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-32 down to 4-33.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value2 = 5
$ws.Range("B3").Value2 = "Macroferia Regional de Talca"
$ws.Range("C3").Value2 = "Maule"
$ws.Range("D3").Value2 = 44756
$ws.Range("E3").Value2 = 7
$ws.Range("F3").Value2 = 100112040
$ws.Range("G3").Value2 = "Cilantro"
$ws.Range("H3").Value2 = "Sin especificar"
$ws.Range("I3").Value2 = "Primera"
$ws.Range("J3").Value2 = 150
$ws.Range("K3").Value2 = 13000
$ws.Range("L3").Value2 = 13000
$ws.Range("M3").Value2 = 13000
$ws.Range("N3").Value2 = "`$/caja 36 atados"
$ws.Range("O3").Value2 = "Región Metropolitana"
$ws.Range("P3").Value2 = 361
$ws.Range("Q3").Value2 = 36
$ws.Range("R3").Value2 = "Hortaliza"
